$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1713.5834
$ws.Range("I40").Value = 1372.2
$ws.Range("J40").Value = 1957.4286
$ws.Range("K40").Value = 1372.2
$ws.Range("L40").Value = 1957.4286
$ws.Range("M40").Value = -1197.2
$ws.Range("N40").Value = -2307.4286
$ws.Range("H51").Value = 3099.9092
$ws.Range("I51").Value = 2799.8572
$ws.Range("J51").Value = 3625
$ws.Range("K51").Value = 2799.8572
$ws.Range("L51").Value = 3625
$ws.Range("M51").Value = -2315.8572
$ws.Range("N51").Value = -4593
$ws.Range("H113").Value = 1733.3334
$ws.Range("I113").Value = 1475
$ws.Range("K113").Value = 1475
$ws.Range("M113").Value = 1779
$ws.Range("H116").Value = 3370.238
$ws.Range("I116").Value = 3857.353
$ws.Range("J116").Value = 1300
$ws.Range("K116").Value = 3857.353
$ws.Range("L116").Value = 1300
$ws.Range("M116").Value = -415.3530000000001
$ws.Range("N116").Value = -8184
$ws.Range("H125").Value = 1729.5625
$ws.Range("I125").Value = 5216
$ws.Range("K125").Value = 46944
$ws.Range("M125").Value = -44484
$ws.Range("H132").Value = 5756.919
$ws.Range("I132").Value = 3589.8076
$ws.Range("J132").Value = 10879.182
$ws.Range("K132").Value = 10769.4228
$ws.Range("L132").Value = 32637.546
$ws.Range("M132").Value = -8239.4228
$ws.Range("N132").Value = -37697.546
$ws.Range("H137").Value = 12502032
$ws.Range("I137").Value = 1599.8889
$ws.Range("J137").Value = 28574016
$ws.Range("K137").Value = 4799.6667
$ws.Range("L137").Value = 85722048
$ws.Range("M137").Value = -2249.6667
$ws.Range("N137").Value = -85727148
$ws.Range("H138").Value = 4718736.5
$ws.Range("I138").Value = 1333.7916
$ws.Range("J138").Value = 8622794
$ws.Range("K138").Value = 4001.3748
$ws.Range("L138").Value = 25868382
$ws.Range("M138").Value = 1138.6252
$ws.Range("N138").Value = -25878662
$ws.Range("H139").Value = 37940
$ws.Range("J139").Value = 37940
$ws.Range("L139").Value = 37940
$ws.Range("N139").Value = -48220
$ws.Range("H140").Value = 74200
$ws.Range("J140").Value = 74200
$ws.Range("L140").Value = 74200
$ws.Range("N140").Value = -84560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11671.244
$ws.Range("I32").Value = 13361.406
$ws.Range("J32").Value = 5661.778
$ws.Range("K32").Value = 13361.406
$ws.Range("L32").Value = 5661.778
$ws.Range("M32").Value = -13074.406
$ws.Range("N32").Value = -6235.778
$ws.Range("H61").Value = 22731658
$ws.Range("I61").Value = 31254364
$ws.Range("J61").Value = 4441.6665
$ws.Range("K61").Value = 31254364
$ws.Range("L61").Value = 4441.6665
$ws.Range("M61").Value = -31254152
$ws.Range("N61").Value = -4865.6665
$ws.Range("H64").Value = 12000
$ws.Range("J64").Value = 12000
$ws.Range("L64").Value = 12000
$ws.Range("N64").Value = -12496
$ws.Range("H67").Value = 12000
$ws.Range("J67").Value = 12000
$ws.Range("L67").Value = 12000
$ws.Range("N67").Value = -13716
$ws.Range("H68").Value = 38099
$ws.Range("J68").Value = 38099
$ws.Range("L68").Value = 38099
$ws.Range("N68").Value = -39721
$ws.Range("H71").Value = 38099
$ws.Range("J71").Value = 38099
$ws.Range("L71").Value = 114297
$ws.Range("N71").Value = -122409
$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25676
$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27340
$ws.Range("H136").Value = 22731658
$ws.Range("I136").Value = 31254364
$ws.Range("J136").Value = 4441.6665
$ws.Range("K136").Value = 93763092
$ws.Range("L136").Value = 13324.9995
$ws.Range("M136").Value = -93760542
$ws.Range("N136").Value = -18424.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 45181
$ws.Range("J62").Value = 45181
$ws.Range("L62").Value = 45181
$ws.Range("N62").Value = -46553
$ws.Range("H65").Value = 45181
$ws.Range("J65").Value = 45181
$ws.Range("L65").Value = 135543
$ws.Range("N65").Value = -142407
$ws.Range("H107").Value = 1553.2142
$ws.Range("I107").Value = 1815.5883
$ws.Range("J107").Value = 1147.7273
$ws.Range("K107").Value = 1815.5883
$ws.Range("L107").Value = 1147.7273
$ws.Range("M107").Value = 104.4117000000001
$ws.Range("N107").Value = -4987.7273
$ws.Range("H134").Value = 4920.3335
$ws.Range("I134").Value = 3253
$ws.Range("J134").Value = 6254.2
$ws.Range("K134").Value = 9759
$ws.Range("L134").Value = 18762.6
$ws.Range("M134").Value = -7224
$ws.Range("N134").Value = -23832.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4046.9092
$ws.Range("I94").Value = 2510.2
$ws.Range("J94").Value = 5327.5
$ws.Range("K94").Value = 2510.2
$ws.Range("L94").Value = 5327.5
$ws.Range("M94").Value = -2059.2
$ws.Range("N94").Value = -6229.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 88.27273
$ws.Range("I23").Value = 60.142857
$ws.Range("J23").Value = 137.5
$ws.Range("K23").Value = 180.428571
$ws.Range("L23").Value = 412.5
$ws.Range("M23").Value = 54.57142899999999
$ws.Range("N23").Value = -882.5
$ws.Range("H86").Value = 1242.5
$ws.Range("I86").Value = 874.5454999999999
$ws.Range("J86").Value = 1692.2222
$ws.Range("K86").Value = 2623.6365
$ws.Range("L86").Value = 5076.6666
$ws.Range("M86").Value = -1437.6365
$ws.Range("N86").Value = -7448.6666
$ws.Range("H89").Value = 1242.5
$ws.Range("I89").Value = 874.5454999999999
$ws.Range("J89").Value = 1692.2222
$ws.Range("K89").Value = 7870.9095
$ws.Range("L89").Value = 15229.9998
$ws.Range("M89").Value = -1942.9095
$ws.Range("N89").Value = -27085.9998
$ws.Range("H104").Value = 10249.5
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H122").Value = 940.3182
$ws.Range("I122").Value = 874.4
$ws.Range("K122").Value = 7869.599999999999
$ws.Range("M122").Value = -5419.599999999999
$ws.Range("H131").Value = 921.8461
$ws.Range("I131").Value = 648.3333
$ws.Range("J131").Value = 1003.9
$ws.Range("K131").Value = 1944.9999
$ws.Range("L131").Value = 3011.7
$ws.Range("M131").Value = 3095.0001
$ws.Range("N131").Value = -13091.7
$ws.Range("H132").Value = 950.4
$ws.Range("I132").Value = 584
$ws.Range("K132").Value = 5256
$ws.Range("M132").Value = -2726
$ws.Range("H136").Value = 3059.5833
$ws.Range("I136").Value = 1539.091
$ws.Range("J136").Value = 4346.154
$ws.Range("K136").Value = 4617.272999999999
$ws.Range("L136").Value = 13038.462
$ws.Range("M136").Value = 482.7270000000008
$ws.Range("N136").Value = -23238.462
$ws.Range("H137").Value = 4594
$ws.Range("I137").Value = 3713.9
$ws.Range("J137").Value = 5111.706
$ws.Range("K137").Value = 11141.7
$ws.Range("L137").Value = 15335.118
$ws.Range("M137").Value = -6041.700000000001
$ws.Range("N137").Value = -25535.118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18400.303
$ws.Range("I70").Value = 24370.166
$ws.Range("J70").Value = 4623.6924
$ws.Range("K70").Value = 24370.166
$ws.Range("L70").Value = 4623.6924
$ws.Range("M70").Value = -24100.166
$ws.Range("N70").Value = -5163.6924
$ws.Range("H73").Value = 18400.303
$ws.Range("I73").Value = 24370.166
$ws.Range("J73").Value = 4623.6924
$ws.Range("K73").Value = 24370.166
$ws.Range("L73").Value = 4623.6924
$ws.Range("M73").Value = -23434.166
$ws.Range("N73").Value = -6495.6924
$ws.Range("H102").Value = 3258.5264
$ws.Range("I102").Value = 3771.4666
$ws.Range("K102").Value = 3771.4666
$ws.Range("M102").Value = -2149.4666
$ws.Range("H126").Value = 3723.8
$ws.Range("I126").Value = 2379.7144
$ws.Range("K126").Value = 7139.1432
$ws.Range("M126").Value = -4669.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4632.9585
$ws.Range("I40").Value = 5815.25
$ws.Range("J40").Value = 3450.6667
$ws.Range("K40").Value = 5815.25
$ws.Range("L40").Value = 3450.6667
$ws.Range("M40").Value = -5679.25
$ws.Range("N40").Value = -3722.6667
$ws.Range("H55").Value = 414.57144
$ws.Range("I55").Value = 239.2
$ws.Range("J55").Value = 574
$ws.Range("K55").Value = 239.2
$ws.Range("L55").Value = 574
$ws.Range("M55").Value = -66.19999999999999
$ws.Range("N55").Value = -920
$ws.Range("H132").Value = 26330074
$ws.Range("I132").Value = 11166.333
$ws.Range("J132").Value = 38477264
$ws.Range("K132").Value = 33498.999
$ws.Range("L132").Value = 115431792
$ws.Range("M132").Value = -30968.999
$ws.Range("N132").Value = -115436852
$ws.Range("H139").Value = 39777.11
$ws.Range("J139").Value = 39668
$ws.Range("L139").Value = 39668
$ws.Range("N139").Value = -49948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 8334333.5
$ws.Range("I3").Value = 12500000
$ws.Range("K3").Value = 12500000
$ws.Range("M3").Value = -12499886
$ws.Range("H74").Value = 6689.8
$ws.Range("J74").Value = 6981.25
$ws.Range("L74").Value = 6981.25
$ws.Range("N74").Value = -8853.25
$ws.Range("H77").Value = 6689.8
$ws.Range("J77").Value = 6981.25
$ws.Range("L77").Value = 20943.75
$ws.Range("N77").Value = -30303.75
$ws.Range("H107").Value = 627.8461
$ws.Range("I107").Value = 692
$ws.Range("J107").Value = 275
$ws.Range("K107").Value = 2076
$ws.Range("L107").Value = 825
$ws.Range("M107").Value = -156
$ws.Range("N107").Value = -4665
$ws.Range("H122").Value = 2586
$ws.Range("I122").Value = 2426.6843
$ws.Range("K122").Value = 7280.0529
$ws.Range("M122").Value = -4830.0529
$ws.Range("H126").Value = 3125.1
$ws.Range("I126").Value = 1250.1111
$ws.Range("K126").Value = 3750.3333
$ws.Range("M126").Value = -1280.3333
$ws.Range("H136").Value = 846
$ws.Range("I136").Value = 846
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2538
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 12
$ws.Range("N136").ClearContents()
